$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new trade record as row 6, matching the existing table layout:
# Date | Profitable | Principle | Start Principle | BuyPrice | SellPrice | IsShortSell | Price Change % | Strong trade
$row = 6

$ws.Cells.Item($row, 1).Value = 42636.589108796295
$ws.Cells.Item($row, 2).Value = $true
$ws.Cells.Item($row, 3).Value = 9962.2900000000009
$ws.Cells.Item($row, 4).Value = 9931.5
$ws.Cells.Item($row, 5).Value = 19.29
$ws.Cells.Item($row, 6).Value = 19.41
$ws.Cells.Item($row, 7).Value = $false
$ws.Cells.Item($row, 8).Value = 0.62
$ws.Cells.Item($row, 9).Value = $false

# Carry over the date/time cell formatting used by column A and G in the
# previous row, instead of assigning a NumberFormat string (which would
# register a brand new custom format code).
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("G5").Copy()
$ws.Range("G6").PasteSpecial(-4122)

$excel.CutCopyMode = $false
